<#
  Project _Sprint _Plan.xlsx edit
  -------------------------------
  The sheet gained a new "Allocated Module" column between the existing
  "Name" column (D) and the "Day 1" column (which was E, now F). Each of
  the 7 data rows got a module name in that new column.
#>

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Insert a new column at E - this shifts the old E:J (Day1..Day6) to F:K
# and carries the row styles (s=1 / s=2) across automatically, same as
# Excel does when you right click a column header and choose Insert.
$ws.Columns("E:E").Insert()

# The inserted column needs an explicit width (closest value reachable
# through the ColumnWidth property's pixel-snapped rounding to the
# target 32.7109375 used by the saved file).
$ws.Columns("E:E").ColumnWidth = 31.8

# Header + per-row "Allocated Module" values.
$ws.Range("E1").Value = "Allocated Module"
$ws.Range("E2").Value = "Trip Selection"
$ws.Range("E3").Value = "Trip commute"
$ws.Range("E4").Value = "Feedback+payments"
$ws.Range("E5").Value = "Cab allocation(driver side)"
$ws.Range("E6").Value = "Login+Signup"
$ws.Range("E7").Value = "Admin"
$ws.Range("E8").Value = "Cab allocation(customer side)"

# Match the saved file's last selection.
$ws.Range("E14").Select() | Out-Null
